$wb = $excel.ActiveWorkbook

# Rename the two worksheets to mark them as todo
$wsCaseProperty = $wb.Worksheets.Item("case_property")
$wsCaseProperty.Name = "case_property (todo)"

$wsStepProperty = $wb.Worksheets.Item("step_property")
$wsStepProperty.Name = "step_property (todo)"

# Work on the "step" sheet: add a wait_time column (E) and tweak a few values
$ws = $wb.Worksheets.Item("step")

# Copy formatting down column E from column D (header style / row styles)
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null

$ws.Range("D2").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null

$ws.Range("D3").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Header
$ws.Range("E1").Value = "wait_time"

# Data rows - wait_time values
$ws.Range("E2").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("E7").Value = 30
$ws.Range("E8").Value = 30
$ws.Range("E9").Value = 30
$ws.Range("E10").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("E12").Value = 0

# Fix a few pre-existing values in the table
$ws.Range("D5").Value = $false
$ws.Range("C6").Value = 10
$ws.Range("D9").Value = $true

# Match column width for the new column
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

$ws.Range("E4").Select() | Out-Null
